# bulk_booking_upload.xlsx — "update bulk upload files"
#
# The sample booking rows on Sheet1 carry hard-coded dates that had
# drifted into the past relative to the e2e fixtures; bump them forward
# by exactly 10 years (the underlying serials only, formatting is left
# untouched) and leave the sheet's selection parked on E4 (the last
# edited cell), matching where the author's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 booking start/end (2017-12-12 -> 2027-12-12)
$ws.Range("D3").Value = 46733.416666666664
$ws.Range("E3").Value = 46733.458333333336

# Row 4 booking start/end (2018-05-23 -> 2028-05-23)
$ws.Range("D4").Value = 46896.5
$ws.Range("E4").Value = 46896.583333333336

# Leave the cursor on E4, the last cell touched.
[void]$ws.Range("E4").Select()
